# CE-QUAL-W2 2025 Workshop Agenda: swap instructor assignments for a
# handful of specific sessions (Hailie <-> Todd) in the two daily
# schedule tables.
#
# The words "Hailie" and "Todd" each appear many times throughout the
# document (in the cover/contact info as well as the Instructor column
# of both schedule tables), so a blind global Find/Replace would touch
# the wrong occurrences. Instead we target the exact table cells that
# need to change, identified by their row's Module/topic text.

$d = $word.ActiveDocument

function Set-InstructorCell($tableIndex, $topicText, $newInstructor) {
    $table = $d.Tables.Item($tableIndex)
    $rowCount = $table.Rows.Count
    for ($r = 1; $r -le $rowCount; $r++) {
        $table = $d.Tables.Item($tableIndex)
        $topicCell = $table.Cell($r, 3)
        $topic = $topicCell.Range.Text.TrimEnd([char]13, [char]7)
        if ($topic -eq $topicText) {
            $table = $d.Tables.Item($tableIndex)
            $instructorCell = $table.Cell($r, 4)
            $instructorCell.Range.Text = $newInstructor
            return
        }
    }
    throw "Row with topic '$topicText' not found in table $tableIndex"
}

# Day 1 table (Tables.Item(1))
Set-InstructorCell 1 "1.04 Lecture - Water Temperature Modeling" "Todd"
Set-InstructorCell 1 "1.05 Lecture - Water Quality Modeling" "Hailie"

# Day 2 table (Tables.Item(2))
Set-InstructorCell 2 "2.03 Lecture - Water Temperature" "Todd"
Set-InstructorCell 2 "2.04 Workshop - Water Temperature" "Todd"
